# Append 7 new order-line rows (rows 7-13) to the Webstaurant Bakery order
# sheet, matching the existing inline-string "text everywhere" layout used
# by rows 1-6 (SKU, Name, Quantity, Cost Per, Total Cost).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("711SPRNKLEPK", "Sprinkles - Pink",          "2", "24.99",  "49.98"),
    @("711SPRNKLEOR", "Sprinkles - Orange",         "1", "25.62",  "25.62"),
    @("3639225768",   "Compound - Raspberry",       "1", "57.49",  "57.49"),
    @("10201311",     "Spice - Italian Seasoning",  "2", "87.89",  "175.78"),
    @("245663CB",     "Box Cake - 6x6x3",           "2", "39.53",  "79.06"),
    @("245885CB",     "Box Cake - 8x8x5",           "1", "34.81",  "34.81"),
    @("409ML90266",   "Choc Curls - Dark",          "4", "127.96", "511.84")
)

$startRow = 7
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    for ($col = 1; $col -le 5; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        # Force text storage so numeric-looking values ("2", "24.99", ...)
        # are written verbatim instead of being coerced into floating point
        # numbers (which would introduce binary rounding noise).
        $cell.NumberFormat = "@"
        $cell.Value = $data[$i][$col - 1]
        # Drop the temporary text format again so the new cells end up with
        # the same (default) styling as the rest of the sheet.
        $cell.ClearFormats()
    }
}
